$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.535604
$ws.Range("N2").Value = 49.606812
$ws.Range("O2").Value = 0.2120453146491552
$ws.Range("P2").Value = 0.2120453146491552
$ws.Range("Q2").Value = 361.299122163608
$ws.Range("R2").Value = 3251.692099472472
$ws.Range("S2").Value = 0.01064471485261256
$ws.Range("T2").Value = 0.01064471485261256

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 40.62063066666667
$ws.Range("N3").Value = 121.861892
$ws.Range("O3").Value = 0.5209011059384622
$ws.Range("P3").Value = 0.5209011059384622
$ws.Range("Q3").Value = 887.5513831607725
$ws.Range("R3").Value = 7987.962448446953
$ws.Range("S3").Value = 0.02614933392091126
$ws.Range("T3").Value = 0.02614933392091126

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 20.825229
$ws.Range("N4").Value = 62.475687
$ws.Range("O4").Value = 0.2670535794123827
$ws.Range("P4").Value = 0.2670535794123827
$ws.Range("Q4").Value = 455.026436080358
$ws.Range("R4").Value = 4095.237924723222
$ws.Range("S4").Value = 0.01340614013527161
$ws.Range("T4").Value = 0.01340614013527161

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 16.535604
$ws.Range("N5").Value = 49.606812
$ws.Range("O5").Value = 0.2120453146491552
$ws.Range("P5").Value = 0.2120453146491552
$ws.Range("Q5").Value = 6367.075168653748
$ws.Range("R5").Value = 57303.67651788374
$ws.Range("S5").Value = 0.1875888853800705
$ws.Range("T5").Value = 0.1875888853800706

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 40.62063066666667
$ws.Range("N6").Value = 121.861892
$ws.Range("O6").Value = 0.5209011059384622
$ws.Range("P6").Value = 0.5209011059384622
$ws.Range("Q6").Value = 15641.07418469796
$ws.Range("R6").Value = 140769.6676622816
$ws.Range("S6").Value = 0.4608225275711437
$ws.Range("T6").Value = 0.4608225275711437

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.825229
$ws.Range("N7").Value = 62.475687
$ws.Range("O7").Value = 0.2670535794123827
$ws.Range("P7").Value = 0.2670535794123827
$ws.Range("Q7").Value = 8018.805871707374
$ws.Range("R7").Value = 72169.25284536637
$ws.Range("S7").Value = 0.236252724478327
$ws.Range("T7").Value = 0.236252724478327

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 16.535604
$ws.Range("N8").Value = 49.606812
$ws.Range("O8").Value = 0.2120453146491552
$ws.Range("P8").Value = 0.2120453146491552
$ws.Range("Q8").Value = 468.7922939543159
$ws.Range("R8").Value = 4219.130645588843
$ws.Range("S8").Value = 0.0138117144164721
$ws.Range("T8").Value = 0.0138117144164721

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 40.62063066666667
$ws.Range("N9").Value = 121.861892
$ws.Range("O9").Value = 0.5209011059384622
$ws.Range("P9").Value = 0.5209011059384622
$ws.Range("Q9").Value = 1151.614336682089
$ws.Range("R9").Value = 10364.5290301388
$ws.Range("S9").Value = 0.03392924444640721
$ws.Range("T9").Value = 0.03392924444640721

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 20.825229
$ws.Range("N10").Value = 62.475687
$ws.Range("O10").Value = 0.2670535794123827
$ws.Range("P10").Value = 0.2670535794123827
$ws.Range("Q10").Value = 590.4052174346909
$ws.Range("R10").Value = 5313.646956912218
$ws.Range("S10").Value = 0.01739471479878406
$ws.Range("T10").Value = 0.01739471479878406
